$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows starting at row 3, pushing the old row 3 data down to row 7
$ws.Rows.Item(3).Resize(4).Insert()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.477828666666666
$ws.Range("H2").Value = 13.433486
$ws.Range("I2").Value = 0.4652827882180238
$ws.Range("J2").Value = 0.4652827882180238
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.08849299999999999
$ws.Range("N2").Value = 0.265479
$ws.Range("O2").Value = 0.9052313210944106
$ws.Range("P2").Value = 0.9052313210944106
$ws.Range("Q2").Value = 0.3962564921993333
$ws.Range("R2").Value = 3.566308429794
$ws.Range("S2").Value = 0.4211885530610925
$ws.Range("T2").Value = 0.4211885530610925

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.477828666666666
$ws.Range("H3").Value = 13.433486
$ws.Range("I3").Value = 0.4652827882180238
$ws.Range("J3").Value = 0.4652827882180238
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.009264333333333334
$ws.Range("N3").Value = 0.027793
$ws.Range("O3").Value = 0.09476867890558938
$ws.Range("P3").Value = 0.09476867890558936
$ws.Range("Q3").Value = 0.04148409737755556
$ws.Range("R3").Value = 0.373356876398
$ws.Range("S3").Value = 0.04409423515693124
$ws.Range("T3").Value = 0.04409423515693123

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.528563666666667
$ws.Range("H4").Value = 7.585691000000001
$ws.Range("I4").Value = 0.2627383137214249
$ws.Range("J4").Value = 0.2627383137214249
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.08849299999999999
$ws.Range("N4").Value = 0.265479
$ws.Range("O4").Value = 0.9052313210944106
$ws.Range("P4").Value = 0.9052313210944106
$ws.Range("Q4").Value = 0.2237601845543333
$ws.Range("R4").Value = 2.013841660989
$ws.Range("S4").Value = 0.2378389508321632
$ws.Range("T4").Value = 0.2378389508321631

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.528563666666667
$ws.Range("H5").Value = 7.585691000000001
$ws.Range("I5").Value = 0.2627383137214249
$ws.Range("J5").Value = 0.2627383137214249
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009264333333333334
$ws.Range("N5").Value = 0.027793
$ws.Range("O5").Value = 0.09476867890558938
$ws.Range("P5").Value = 0.09476867890558936
$ws.Range("Q5").Value = 0.02342545666255556
$ws.Range("R5").Value = 0.210829109963
$ws.Range("S5").Value = 0.02489936288926172
$ws.Range("T5").Value = 0.02489936288926172

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.617494
$ws.Range("H6").Value = 7.852482
$ws.Range("I6").Value = 0.2719788980605514
$ws.Range("J6").Value = 0.2719788980605514
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.08849299999999999
$ws.Range("N6").Value = 0.265479
$ws.Range("O6").Value = 0.9052313210944106
$ws.Range("P6").Value = 0.9052313210944106
$ws.Range("Q6").Value = 0.231629896542
$ws.Range("R6").Value = 2.084669068878
$ws.Range("S6").Value = 0.2462038172011549
$ws.Range("T6").Value = 0.2462038172011549

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.617494
$ws.Range("H7").Value = 7.852482
$ws.Range("I7").Value = 0.2719788980605514
$ws.Range("J7").Value = 0.2719788980605514
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.009264333333333334
$ws.Range("N7").Value = 0.027793
$ws.Range("O7").Value = 0.09476867890558938
$ws.Range("P7").Value = 0.09476867890558936
$ws.Range("Q7").Value = 0.024249336914
$ws.Range("R7").Value = 0.218244032226
$ws.Range("S7").Value = 0.02577508085939642
$ws.Range("T7").Value = 0.02577508085939641
